$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3, shifting all existing rows
# (old row 3 -> new row 4, ..., old row 101 -> new row 102) down by one.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new weekly record.
$ws.Cells.Item(3, 1).Value = 6
$ws.Cells.Item(3, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(3, 3).Value = "Metropolitana"
$ws.Cells.Item(3, 4).Value = 44812
$ws.Cells.Item(3, 5).Value = 13
$ws.Cells.Item(3, 6).Value = 100114007
$ws.Cells.Item(3, 7).Value = "Jengibre"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 300
$ws.Cells.Item(3, 11).Value = 10000
$ws.Cells.Item(3, 12).Value = 11000
$ws.Cells.Item(3, 13).Value = 10500
$ws.Cells.Item(3, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(3, 15).Value = "Perú"
$ws.Cells.Item(3, 16).Value = 808
$ws.Cells.Item(3, 17).Value = 13
$ws.Cells.Item(3, 18).Value = "Hortaliza"
